{"js": "// Replace each three-digit-by-one-digit multiplication expression\n// in the document with its updated value, matching the exact\n// old text (the expressions are unique, so a plain search/replace\n// on each pair is safe and will not double-replace).\nconst replacements = [\n  [\"918\u00d77=\", \"867\u00d74=\"],\n  [\"469\u00d74=\", \"817\u00d77=\"],\n  [\"184\u00d74=\", \"688\u00d73=\"],\n  [\"330\u00d78=\", \"937\u00d77=\"],\n  [\"666\u00d72=\", \"755\u00d72=\"],\n  [\"238\u00d72=\", \"286\u00d77=\"],\n  [\"824\u00d79=\", \"160\u00d74=\"],\n  [\"232\u00d73=\", \"376\u00d76=\"],\n  [\"575\u00d79=\", \"178\u00d79=\"],\n  [\"205\u00d78=\", \"441\u00d77=\"],\n  [\"524\u00d78=\", \"445\u00d73=\"],\n  [\"473\u00d79=\", \"535\u00d72=\"],\n  [\"742\u00d74=\", \"738\u00d74=\"],\n  [\"360\u00d75=\", \"121\u00d78=\"],\n  [\"683\u00d78=\", \"734\u00d72=\"],\n  [\"751\u00d76=\", \"709\u00d79=\"],\n  [\"369\u00d76=\", \"297\u00d73=\"],\n  [\"303\u00d78=\", \"691\u00d76=\"],\n  [\"423\u00d75=\", \"235\u00d78=\"],\n  [\"454\u00d76=\", \"622\u00d73=\"],\n  [\"924\u00d72=\", \"417\u00d79=\"],\n  [\"469\u00d79=\", \"976\u00d74=\"],\n  [\"399\u00d72=\", \"500\u00d76=\"],\n  [\"913\u00d73=\", \"357\u00d72=\"],\n  [\"987\u00d77=\", \"449\u00d77=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update each three-digit-by-one-digit multiplication expression\n# in the document with its new value. Each old expression is\n# unique in the document, so Find/Replace with MatchCase and\n# ReplaceAll is safe (exactly one hit per pair, no cross-matches).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"918\u00d77=\", \"867\u00d74=\"),\n    @(\"469\u00d74=\", \"817\u00d77=\"),\n    @(\"184\u00d74=\", \"688\u00d73=\"),\n    @(\"330\u00d78=\", \"937\u00d77=\"),\n    @(\"666\u00d72=\", \"755\u00d72=\"),\n    @(\"238\u00d72=\", \"286\u00d77=\"),\n    @(\"824\u00d79=\", \"160\u00d74=\"),\n    @(\"232\u00d73=\", \"376\u00d76=\"),\n    @(\"575\u00d79=\", \"178\u00d79=\"),\n    @(\"205\u00d78=\", \"441\u00d77=\"),\n    @(\"524\u00d78=\", \"445\u00d73=\"),\n    @(\"473\u00d79=\", \"535\u00d72=\"),\n    @(\"742\u00d74=\", \"738\u00d74=\"),\n    @(\"360\u00d75=\", \"121\u00d78=\"),\n    @(\"683\u00d78=\", \"734\u00d72=\"),\n    @(\"751\u00d76=\", \"709\u00d79=\"),\n    @(\"369\u00d76=\", \"297\u00d73=\"),\n    @(\"303\u00d78=\", \"691\u00d76=\"),\n    @(\"423\u00d75=\", \"235\u00d78=\"),\n    @(\"454\u00d76=\", \"622\u00d73=\"),\n    @(\"924\u00d72=\", \"417\u00d79=\"),\n    @(\"469\u00d79=\", \"976\u00d74=\"),\n    @(\"399\u00d72=\", \"500\u00d76=\"),\n    @(\"913\u00d73=\", \"357\u00d72=\"),\n    @(\"987\u00d77=\", \"449\u00d77=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $oldText,    # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $newText,    # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n"}
